# Weekly refresh of the "Hortaliza, Macroferia Regional de Talca - Berenjena"
# sheet: a new week's record is inserted at row 15 (pushing all subsequent
# daily records down by one row), which grows the used range from A1:R61
# to A1:R62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; Excel shifts rows 15-61 down to 16-62
# and copies formatting (incl. the date-time number format on column D)
# from the row above, just like an interactive "Insert Row" would.
$ws.Rows.Item(15).Insert()

# Populate the newly-inserted row with the new week's record.
$ws.Range("A15").Value = 5
$ws.Range("B15").Value = "Macroferia Regional de Talca"
$ws.Range("C15").Value = "Maule"
$ws.Range("D15").Value = "9/20/2021"
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = 100112001
$ws.Range("G15").Value = "Berenjena"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 7000
$ws.Range("N15").Value = "$/caja 60 unidades"
$ws.Range("O15").Value = "Región de Arica y Parinacota"
$ws.Range("P15").Value = 117
$ws.Range("Q15").Value = 60
$ws.Range("R15").Value = "Hortaliza"
